$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column P (29-jun) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("O1").Copy()
$ws1.Range("P1").PasteSpecial(-4122)
$ws1.Range("P1").Value = "29-jun"

$p2_25 = @(95.06999999999999, 89.42, 89.66, 86.81999999999999, 84.3, 85.3, 86, 74.70999999999999, 47.4, 3.75, 0, -0.08, -0.02, -0.03, -0.09, -0.01, 0, 20, 79.38, 102.14, 111.98, 118.07, 113.73, 101)

for ($i = 0; $i -lt $p2_25.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 16).Value = $p2_25[$i]
}

# --- Sheet "Gaz": add row 13 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A13").NumberFormat = "@"
$ws2.Range("A13").Value = "2025-06-27"
$ws2.Range("A13").Style = "Normal"
$ws2.Range("B13").Value = 32.7

# --- Sheet "CO2": add row 13 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A13").NumberFormat = "@"
$ws3.Range("A13").Value = "2025-06-27"
$ws3.Range("A13").Style = "Normal"
$ws3.Range("B13").Value = 69.92
